$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "ticket" placeholder namespace to "tickets" across the two
# templated rows (row 2 = d.ticket[i].*, row 3 = d.ticket[i + 1].*).
for ($r = 2; $r -le 3; $r++) {
  for ($c = 1; $c -le 7; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $val = $cell.Value()
    if ($val -ne $null) {
      $newVal = $val.Replace("{d.ticket[", "{d.tickets[")
      $cell.Value = $newVal
    }
  }
}

# Move the active selection on the sheet.
$ws.Range("F19").Select()
